$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# NOTE: Price/Volume/Hora columns store numeric-looking values as literal TEXT
# in the source data (t="inlineStr"), not as numbers/percentages. Force the
# target cells to text format before assigning so Excel does not silently
# reinterpret them as numbers (which would also strip formatting like "1.70%").

# --- Column B (Coin) changes (rows shifted: HotbitToken now sorts first) ---
$ws.Range("B19").Value = 'HotbitToken'
$ws.Range("B20").Value = 'LEO'
$ws.Range("B21").Value = 'BitpandaEcosystemToken'
$ws.Range("B22").Value = 'MCDex'
$ws.Range("B23").Value = 'ProBitToken'
$ws.Range("B24").Value = 'ZBToken'
$ws.Range("B25").Value = 'CoinExToken'
$ws.Range("B26").Value = 'BitKan'

# --- Column C (Link) changes ---
$ws.Range("C19").Value = 'https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb'
$ws.Range("C20").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("C21").Value = 'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best'
$ws.Range("C22").Value = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
$ws.Range("C23").Value = 'https://coinranking.com/coin/lQP4d6T2+probittoken-prob'
$ws.Range("C24").Value = 'https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb'
$ws.Range("C25").Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
$ws.Range("C26").Value = 'https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan'

# --- Column D (Price) changes ---
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '328.68'
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '41.52'
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '5.624'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.08171'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '2.026'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '8.733'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '4.525'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.949'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.9174'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.1276'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.1951'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.09445'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.03802'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.001298'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.006225'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.004405'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '3.440'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.3494'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '8.248'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.1394'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.2412'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.04417'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.001259'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.02762'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.05415'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.007661'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.1415'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.009006'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.002132'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.00006607'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.003241'

# --- Column E (Volume(1h)) changes ---
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '1.70%'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '4.85%'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '-3.93%'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '1.93%'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '-3.09%'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '1.07%'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '-1.06%'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '0.21%'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '-1.36%'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '0.10%'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '-0.36%'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '3.20%'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '6.28%'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '1.25%'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '0.40%'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '0.06%'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '2.28%'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '2.71%'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '-1.24%'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '-5.66%'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '1.63%'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '-1.58%'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '0.11%'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '-0.23%'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '2.60%'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '10.16%'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '3.03%'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '2.91%'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '0.76%'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '-6.32%'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '0.64%'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '15.74%'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '-2.06%'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '0.08%'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '7.93%'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '0.08%'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '0.08%'

# --- Column G (Hora): every data row advances from 6 to 7 ---
$ws.Range("G2:G51").NumberFormat = "@"
$ws.Range("G2:G51").Value = "7"

